$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = -6.078699999999995
$ws.Range("D27").Value = -7.908400000000003
$ws.Range("D32").Value = -7.310799999999993
$ws.Range("D36").Value = -7.878099999999997
$ws.Range("D38").Value = -7.226899999999998
$ws.Range("D46").Value = -8.106299999999994
$ws.Range("D54").Value = -7.853000000000002
$ws.Range("D55").Value = -7.242899999999996
$ws.Range("D56").Value = -8.9641
$ws.Range("D67").Value = -7.347399999999994
$ws.Range("D69").Value = -7.614999999999997
$ws.Range("D72").Value = -7.378999999999998
$ws.Range("D83").Value = -9.250500000000001
$ws.Range("D86").Value = -8.476400000000002
$ws.Range("D91").Value = -8.016300000000001
$ws.Range("D93").Value = -6.988299999999994
$ws.Range("D99").Value = -8.104800000000004
$ws.Range("D104").Value = -7.642500000000001
